$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E52").Value = " 2024-09-02"
$ws.Range("E54").Value = " 2024-09-25"
$ws.Range("E58").Value = " 2024-10-23"
$ws.Range("E62").Value = " 2024-11-20"
$ws.Range("E63").Value = " 2024-11-28"
$ws.Range("E60").Value = " 2024-11-11"

[void]$ws.Range("E60").Select()
$excel.ActiveWindow.ScrollRow = 24
